$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows above row 382 (current latest weekly block), shifting
# all existing data (old rows 382:442) down to (386:446).
$ws.Rows("382:385").Insert()

# Fill the newly inserted rows with the new weekly snapshot
# (date 2023-07-24 -> serial 45131), variety "Sin especificar",
# qualities Extra / Primera / Segunda / Tercera.

# Row 382 - Extra
$ws.Range("A382").Value = 6
$ws.Range("B382").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C382").Value = "Metropolitana"
$ws.Range("D382").Value = 45131
$ws.Range("E382").Value = 13
$ws.Range("F382").Value = 100112043
$ws.Range("G382").Value = "Pepino dulce"
$ws.Range("H382").Value = "Sin especificar"
$ws.Range("I382").Value = "Extra"
$ws.Range("J382").Value = 270
$ws.Range("K382").Value = 16000
$ws.Range("L382").Value = 16000
$ws.Range("M382").Value = 16000
$ws.Range("N382").Value = "$/bandeja 18 kilos"
$ws.Range("O382").Value = "Provincia de Limarí"
$ws.Range("P382").Value = 889
$ws.Range("Q382").Value = 18
$ws.Range("R382").Value = "Hortaliza"

# Row 383 - Primera
$ws.Range("A383").Value = 6
$ws.Range("B383").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C383").Value = "Metropolitana"
$ws.Range("D383").Value = 45131
$ws.Range("E383").Value = 13
$ws.Range("F383").Value = 100112043
$ws.Range("G383").Value = "Pepino dulce"
$ws.Range("H383").Value = "Sin especificar"
$ws.Range("I383").Value = "Primera"
$ws.Range("J383").Value = 340
$ws.Range("K383").Value = 15000
$ws.Range("L383").Value = 15000
$ws.Range("M383").Value = 15000
$ws.Range("N383").Value = "$/bandeja 18 kilos"
$ws.Range("O383").Value = "Provincia de Limarí"
$ws.Range("P383").Value = 833
$ws.Range("Q383").Value = 18
$ws.Range("R383").Value = "Hortaliza"

# Row 384 - Segunda
$ws.Range("A384").Value = 6
$ws.Range("B384").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C384").Value = "Metropolitana"
$ws.Range("D384").Value = 45131
$ws.Range("E384").Value = 13
$ws.Range("F384").Value = 100112043
$ws.Range("G384").Value = "Pepino dulce"
$ws.Range("H384").Value = "Sin especificar"
$ws.Range("I384").Value = "Segunda"
$ws.Range("J384").Value = 340
$ws.Range("K384").Value = 13000
$ws.Range("L384").Value = 13000
$ws.Range("M384").Value = 13000
$ws.Range("N384").Value = "$/bandeja 18 kilos"
$ws.Range("O384").Value = "Provincia de Limarí"
$ws.Range("P384").Value = 722
$ws.Range("Q384").Value = 18
$ws.Range("R384").Value = "Hortaliza"

# Row 385 - Tercera
$ws.Range("A385").Value = 6
$ws.Range("B385").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C385").Value = "Metropolitana"
$ws.Range("D385").Value = 45131
$ws.Range("E385").Value = 13
$ws.Range("F385").Value = 100112043
$ws.Range("G385").Value = "Pepino dulce"
$ws.Range("H385").Value = "Sin especificar"
$ws.Range("I385").Value = "Tercera"
$ws.Range("J385").Value = 175
$ws.Range("K385").Value = 10000
$ws.Range("L385").Value = 10000
$ws.Range("M385").Value = 10000
$ws.Range("N385").Value = "$/bandeja 18 kilos"
$ws.Range("O385").Value = "Provincia de Limarí"
$ws.Range("P385").Value = 556
$ws.Range("Q385").Value = 18
$ws.Range("R385").Value = "Hortaliza"

# Preserve the date-time number format for the new D-column cells, matching
# the style used throughout column D in this sheet.
$ws.Range("D382:D385").NumberFormat = "YYYY-MM-DD HH:MM:SS"
